$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1929260450160772
$ws.Range("C2").Value = 0.5819935691318328
$ws.Range("J2").Value = 0.009646302250803859
$ws.Range("P2").Value = 0.1543408360128617
$ws.Range("S2").Value = 0.06109324758842444
$ws.Range("B3").Value = 0.0108695652173913
$ws.Range("C3").Value = 0.02173913043478261
$ws.Range("J3").Value = 0.03260869565217391
$ws.Range("P3").Value = 0.7173913043478261
$ws.Range("S3").Value = 0.2173913043478261
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.6904761904761905
$ws.Range("S4").Value = 0.2380952380952381
$ws.Range("B6").Value = 0.04721030042918455
$ws.Range("D6").Value = 0.01716738197424893
$ws.Range("F6").Value = 0.07725321888412018
$ws.Range("J6").Value = 0.2618025751072962
$ws.Range("O6").Value = 0.008583690987124463
$ws.Range("Q6").Value = 0.167381974248927
$ws.Range("R6").Value = 0.09442060085836911
$ws.Range("S6").Value = 0.3261802575107296
$ws.Range("B7").Value = 0.06572769953051644
$ws.Range("D7").Value = 0.0187793427230047
$ws.Range("E7").Value = 0.004694835680751174
$ws.Range("F7").Value = 0.07511737089201878
$ws.Range("J7").Value = 0.1690140845070423
$ws.Range("O7").Value = 0.01408450704225352
$ws.Range("Q7").Value = 0.1596244131455399
$ws.Range("R7").Value = 0.1032863849765258
$ws.Range("S7").Value = 0.3896713615023474
$ws.Range("B8").Value = 0.09045226130653267
$ws.Range("D8").Value = 0.008375209380234505
$ws.Range("F8").Value = 0.06867671691792294
$ws.Range("J8").Value = 0.1155778894472362
$ws.Range("O8").Value = 0.02177554438860971
$ws.Range("Q8").Value = 0.1742043551088777
$ws.Range("R8").Value = 0.09212730318257957
$ws.Range("S8").Value = 0.4288107202680067
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.01304347826086956
$ws.Range("E9").Value = 0.004347826086956522
$ws.Range("F9").Value = 0.05217391304347826
$ws.Range("J9").Value = 0.1173913043478261
$ws.Range("O9").Value = 0.02608695652173913
$ws.Range("Q9").Value = 0.1739130434782609
$ws.Range("R9").Value = 0.08260869565217391
$ws.Range("S9").Value = 0.4304347826086957
$ws.Range("B10").Value = 0.1007142857142857
$ws.Range("D10").Value = 0.01928571428571428
$ws.Range("F10").Value = 0.055
$ws.Range("J10").Value = 0.1278571428571429
$ws.Range("O10").Value = 0.01071428571428571
$ws.Range("Q10").Value = 0.22
$ws.Range("R10").Value = 0.09
$ws.Range("S10").Value = 0.3764285714285714
$ws.Range("G11").Value = 0.1517615176151761
$ws.Range("J11").Value = 0.08672086720867209
$ws.Range("K11").Value = 0.2195121951219512
$ws.Range("L11").Value = 0.4986449864498645
$ws.Range("S11").Value = 0.04336043360433604
$ws.Range("G12").Value = 0.6875
$ws.Range("J12").Value = 0.1770833333333333
$ws.Range("K12").Value = 0.02083333333333333
$ws.Range("L12").Value = 0.05208333333333334
$ws.Range("S12").Value = 0.0625
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2325581395348837
$ws.Range("S13").Value = 0.06976744186046512
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01376146788990826
$ws.Range("H15").Value = 0.1651376146788991
$ws.Range("I15").Value = 0.02752293577981652
$ws.Range("J15").Value = 0.3990825688073394
$ws.Range("K15").Value = 0.06880733944954129
$ws.Range("M15").Value = 0.009174311926605505
$ws.Range("O15").Value = 0.05504587155963303
$ws.Range("S15").Value = 0.2614678899082569
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.2315270935960591
$ws.Range("I16").Value = 0.08866995073891626
$ws.Range("J16").Value = 0.4137931034482759
$ws.Range("K16").Value = 0.09852216748768473
$ws.Range("M16").Value = 0.009852216748768473
$ws.Range("O16").Value = 0.01970443349753695
$ws.Range("S16").Value = 0.1231527093596059
$ws.Range("F17").Value = 0.01919385796545105
$ws.Range("H17").Value = 0.2034548944337812
$ws.Range("I17").Value = 0.08829174664107485
$ws.Range("J17").Value = 0.3608445297504799
$ws.Range("K17").Value = 0.1036468330134357
$ws.Range("M17").Value = 0.01727447216890595
$ws.Range("O17").Value = 0.06525911708253358
$ws.Range("S17").Value = 0.1420345489443378
$ws.Range("F18").Value = 0.00823045267489712
$ws.Range("H18").Value = 0.2098765432098765
$ws.Range("I18").Value = 0.09465020576131687
$ws.Range("J18").Value = 0.3909465020576132
$ws.Range("K18").Value = 0.102880658436214
$ws.Range("M18").Value = 0.00411522633744856
$ws.Range("N18").Value = 0.00411522633744856
$ws.Range("O18").Value = 0.04526748971193416
$ws.Range("S18").Value = 0.139917695473251
$ws.Range("F19").Value = 0.01322751322751323
$ws.Range("H19").Value = 0.2380952380952381
$ws.Range("I19").Value = 0.09193121693121693
$ws.Range("J19").Value = 0.3333333333333333
$ws.Range("K19").Value = 0.1097883597883598
$ws.Range("M19").Value = 0.0205026455026455
$ws.Range("N19").Value = 0.002645502645502645
$ws.Range("O19").Value = 0.05753968253968254
$ws.Range("S19").Value = 0.1329365079365079
